$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Row 2
Set-TextValue $ws.Range('D2') '68.802.72'
Set-TextValue $ws.Range('E2') '  -0.35%  '
# Row 3
Set-TextValue $ws.Range('D3') '3.836.07'
Set-TextValue $ws.Range('E3') '  +2.40%  '
# Row 4
Set-TextValue $ws.Range('E4') '  +0.13%  '
# Row 5
Set-TextValue $ws.Range('D5') '601.63'
Set-TextValue $ws.Range('E5') '  -0.11%  '
# Row 6
Set-TextValue $ws.Range('D6') '162.93'
Set-TextValue $ws.Range('E6') '  -2.92%  '
# Row 7
Set-TextValue $ws.Range('D7') '3.831.68'
Set-TextValue $ws.Range('E7') '  +2.31%  '
# Row 8
Set-TextValue $ws.Range('E8') '  +0.18%  '
# Row 9
Set-TextValue $ws.Range('D9') '0.531'
Set-TextValue $ws.Range('E9') '  -1.77%  '
# Row 10
Set-TextValue $ws.Range('D10') '0.167'
Set-TextValue $ws.Range('E10') '  -0.11%  '
# Row 11
Set-TextValue $ws.Range('D11') '6.30'
Set-TextValue $ws.Range('E11') '  -2.37%  '
# Row 12
Set-TextValue $ws.Range('D12') '0.458'
Set-TextValue $ws.Range('E12') '  -0.41%  '
# Row 13
Set-TextValue $ws.Range('D13') '36.85'
Set-TextValue $ws.Range('E13') '  -3.29%  '
# Row 14
Set-TextValue $ws.Range('D14') '0.0000243'
Set-TextValue $ws.Range('E14') '  -1.68%  '
# Row 15
Set-TextValue $ws.Range('D15') '4.479.53'
Set-TextValue $ws.Range('E15') '  +2.44%  '
# Row 16
Set-TextValue $ws.Range('D16') '3.852.51'
Set-TextValue $ws.Range('E16') '  +2.77%  '
# Row 17
Set-TextValue $ws.Range('D17') '69.018.37'
Set-TextValue $ws.Range('E17') '  -0.06%  '
# Row 18
Set-TextValue $ws.Range('D18') '7.57'
Set-TextValue $ws.Range('E18') '  +3.14%  '
# Row 19
Set-TextValue $ws.Range('D19') '11.53'
Set-TextValue $ws.Range('E19') '  +5.01%  '
# Row 20
Set-TextValue $ws.Range('E20') '  +0.00%  '
# Row 21
Set-TextValue $ws.Range('D21') '17.12'
Set-TextValue $ws.Range('E21') '  -0.52%  '
# Row 22
Set-TextValue $ws.Range('D22') '485.59'
Set-TextValue $ws.Range('E22') '  -1.45%  '
# Row 23
Set-TextValue $ws.Range('D23') '0.718'
Set-TextValue $ws.Range('E23') '  -0.94%  '
# Row 24
Set-TextValue $ws.Range('D24') '0.0000158'
Set-TextValue $ws.Range('E24') '  +4.52%  '
# Row 25
Set-TextValue $ws.Range('D25') '84.36'
Set-TextValue $ws.Range('E25') '  -0.56%  '
# Row 26
Set-TextValue $ws.Range('D26') '2.25'
Set-TextValue $ws.Range('E26') '  -2.66%  '
# Row 27
Set-TextValue $ws.Range('D27') '12.15'
Set-TextValue $ws.Range('E27') '  -1.25%  '
# Row 28
Set-TextValue $ws.Range('D28') '0.999'
Set-TextValue $ws.Range('E28') '  -0.04%  '
# Row 29
Set-TextValue $ws.Range('D29') '9.94'
Set-TextValue $ws.Range('E29') '  -2.05%  '
# Row 30
Set-TextValue $ws.Range('D30') '2.97'
Set-TextValue $ws.Range('E30') '  -0.77%  '
# Row 31
Set-TextValue $ws.Range('D31') '7.91'
Set-TextValue $ws.Range('E31') '  -1.91%  '
# Row 32
Set-TextValue $ws.Range('B32') 'WrappedeETH'
Set-TextValue $ws.Range('C32') 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
Set-TextValue $ws.Range('D32') '4.005.34'
Set-TextValue $ws.Range('E32') '  +2.90%  '
# Row 33
Set-TextValue $ws.Range('B33') 'ImmutableX'
Set-TextValue $ws.Range('C33') 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Range('D33') '2.38'
Set-TextValue $ws.Range('E33') '  -4.09%  '
# Row 34
Set-TextValue $ws.Range('D34') '32.14'
Set-TextValue $ws.Range('E34') '  +1.75%  '
# Row 35
Set-TextValue $ws.Range('D35') '3.783.54'
Set-TextValue $ws.Range('E35') '  +2.73%  '
# Row 36
Set-TextValue $ws.Range('D36') '0.107'
Set-TextValue $ws.Range('E36') '  -1.89%  '
# Row 37
Set-TextValue $ws.Range('D37') '1.03'
Set-TextValue $ws.Range('E37') '  +1.49%  '
# Row 38
Set-TextValue $ws.Range('D38') '0.139'
Set-TextValue $ws.Range('E38') '  +3.74%  '
# Row 39
Set-TextValue $ws.Range('D39') '5.86'
Set-TextValue $ws.Range('E39') '  -0.17%  '
# Row 40
Set-TextValue $ws.Range('E40') '  +0.30%  '
# Row 41
Set-TextValue $ws.Range('D41') '0.318'
Set-TextValue $ws.Range('E41') '  -2.23%  '
# Row 42
Set-TextValue $ws.Range('D42') '2.98'
Set-TextValue $ws.Range('E42') '  +0.47%  '
# Row 43
Set-TextValue $ws.Range('D43') '436.26'
Set-TextValue $ws.Range('E43') '  +0.89%  '
# Row 44
Set-TextValue $ws.Range('B44') 'OKB'
Set-TextValue $ws.Range('C44') 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws.Range('D44') '48.46'
Set-TextValue $ws.Range('E44') '  -0.36%  '
# Row 45
Set-TextValue $ws.Range('B45') 'Stacks'
Set-TextValue $ws.Range('C45') 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Range('D45') '1.98'
Set-TextValue $ws.Range('E45') '  -0.76%  '
# Row 46
Set-TextValue $ws.Range('B46') 'USDe'
Set-TextValue $ws.Range('C46') 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue $ws.Range('D46') '1.00'
Set-TextValue $ws.Range('E46') '  -0.01%  '
# Row 47
Set-TextValue $ws.Range('B47') 'Cosmos'
Set-TextValue $ws.Range('C47') 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue $ws.Range('D47') '8.38'
Set-TextValue $ws.Range('E47') '  -1.65%  '
# Row 48
Set-TextValue $ws.Range('B48') 'EnergySwap'
Set-TextValue $ws.Range('C48') 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range('D48') '27.34'
Set-TextValue $ws.Range('E48') '  +17.29%  '
# Row 49
Set-TextValue $ws.Range('D49') '2.845.16'
Set-TextValue $ws.Range('E49') '  +1.94%  '
# Row 50
Set-TextValue $ws.Range('D50') '142.16'
Set-TextValue $ws.Range('E50') '  +0.47%  '
# Row 51
Set-TextValue $ws.Range('D51') '0.0356'
Set-TextValue $ws.Range('E51') '  +0.82%  '
